# Add a hyperlink to the GitHub URL shown on the "GitHub Link" slide, and
# press Enter after it to leave a trailing blank line/paragraph - matching
# the target edit captured in the diff (hlinkClick on the run, plus a new
# empty paragraph after it).

$p = $ppt.ActivePresentation

$targetUrl = "https://github.com/tanishverse/edunet_project_skillsbuild"

# Locate the slide/shape that holds the GitHub URL instead of hard-coding
# indices, so the script is resilient to the deck being reshuffled.
$slide = $null
$shape = $null
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $candidateSlide = $p.Slides.Item($i)
    for ($j = 1; $j -le $candidateSlide.Shapes.Count; $j++) {
        $candidateShape = $candidateSlide.Shapes.Item($j)
        if ($candidateShape.HasTextFrame) {
            if ($candidateShape.TextFrame.TextRange.Text -eq $targetUrl) {
                $slide = $candidateSlide
                $shape = $candidateShape
            }
        }
    }
}

$tr = $shape.TextFrame.TextRange
$urlLength = $tr.Length

# Press Enter at the end of the URL line - this creates a second, empty
# paragraph right after the URL paragraph.
$newPara = $tr.InsertAfter("`r")

# Re-grab the (now two-paragraph) text range and apply the hyperlink only to
# the original URL run, leaving the new trailing paragraph mark untouched.
$full = $shape.TextFrame.TextRange
$urlRange = $full.Characters(1, $urlLength)
$urlRange.ActionSettings.Item(1).Hyperlink.Address = $targetUrl
